# Update "ECO Actual" (column B) and "ECO Balance" (column D) figures
# for each FSR row. Values in this sheet are stored as text (e.g. "45.00"),
# so force a text number format before writing so the trailing zeros /
# text type are preserved instead of Excel auto-converting to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = "24.00"
    3 = "50.00"
    4 = "58.00"
    5 = "50.00"
    6 = "38.00"
    7 = "220.00"
}

$bRange = $ws.Range("B2:B7")
$dRange = $ws.Range("D2:D7")
$bRange.NumberFormat = "@"
$dRange.NumberFormat = "@"

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Range("B$row").Value = $value
    $ws.Range("D$row").Value = $value
}
